$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header B1 ("sdfsdf" -> "sdfsdf_1")
$ws.Range("B1").Value = "sdfsdf_1"

# Add new header columns F1:J1
$ws.Range("F1").Value = "chatId"
$ws.Range("G1").Value = "firstName"
$ws.Range("H1").Value = "lastName"
$ws.Range("I1").Value = "email"
$ws.Range("J1").Value = "imageUrl"

# The old row 3 user (Huntrick) replaces row 2 (Moin/Looser is dropped)
$ws.Range("A2").Value = "Huntrick"
$ws.Range("B2").Value = "GB"
$ws.Range("C2").Value = "huntrick@gmail.com"
$ws.Range("D2").Value = "USDT"
$ws.Range("E2").Value = "AgACAgIAAxkBAAIBiWfhtJaGFGUSa8oyqUl4J5uzqNwaAAJT7jEbnVMIS1hpkXQDSbihAQADAgADeAADNgQ"

# Row 3 no longer carries the old A:E user columns
$ws.Range("A3:E3").ClearContents()

# Row 3 now holds a new telegram-linked user record in F:J
$ws.Range("F3").Value = 5038824563
$ws.Range("G3").Value = "sadfasdf"
$ws.Range("H3").Value = "dfasdfsad"
$ws.Range("I3").Value = "asdfasdf@sdfasdf.com"
$ws.Range("J3").Value = "https://api.telegram.org/file/bot7203212788:AAF46vpGtp-c0vGf8Twd3flL4tFSFvZKMXM/photos/file_0.jpg"

# Turn off right-to-left sheet view
$excel.ActiveWindow.DisplayRightToLeft = $false
